$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '41.527.92'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.41%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.482.47'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +0.71%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.21'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.08%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '92.76'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.46%  '
$ws.Range('E7').Value = '  -1.61%  '
$ws.Range('E8').Value = '  -0.13%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.494'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -2.60%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '32.29'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -3.85%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0781'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.13%  '
$ws.Range('E12').Value = '  +1.34%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.867.26'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.96%  '
$ws.Range('E14').Value = '  -2.05%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.26'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +4.89%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.472.56'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.27%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.757'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -4.14%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '41.610.57'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.65%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.29'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -1.27%  '
$ws.Range('E20').Value = '  +0.17%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '70.57'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +2.37%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.04'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -4.53%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.23'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -1.12%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.68'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -3.12%  '
$ws.Range('E25').Value = '  -0.29%  '
$ws.Range('E26').Value = '  -2.98%  '
$ws.Range('E27').Value = '  -1.92%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.23'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.75%  '
$ws.Range('E29').Value = '  -1.68%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '35.92'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -1.30%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '153.36'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.06%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.38'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -4.77%  '
$ws.Range('E33').Value = '  -2.20%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0755'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.20%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '17.90'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +4.43%  '
$ws.Range('E36').Value = '  -1.04%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.96'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -2.07%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.82'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -3.10%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.112'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.55%  '
$ws.Range('E40').Value = '  -3.80%  '
$ws.Range('E41').Value = '  -0.42%  '
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '19.46'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -8.75%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.940.70'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -2.44%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0282'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.92%  '
$ws.Range('E46').Value = '  -3.81%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.77'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.06%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.728.13'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.00%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '95.59'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -2.12%  '
$ws.Range('E50').Value = '  -3.11%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '66.78'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -3.81%  '
